{"js": "// Add a new list item \"Bhuratan\" after the last paragraph (\"Tushar\"),\n// matching the existing list paragraph formatting, and move the\n// \"_GoBack\" bookmark so it sits at the end of the new last paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the list (\"Tushar\") - append the new name after it.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"Bhuratan\", \"After\");\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark from around \"Tushar\" to just after \"Bhuratan\".\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst found = body.search(\"Bhuratan\", { matchCase: true, matchWholeWord: true });\nfound.load(\"items\");\nawait context.sync();\n\nconst newNameRange = found.items[0].getRange(\"End\");\nnewNameRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Add a new list item \"Bhuratan\" right after the last paragraph (\"Tushar\"),\n# inheriting the list/paragraph formatting automatically via InsertParagraphAfter.\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Bhuratan\"\n\n# Move the \"_GoBack\" bookmark off of \"Tushar\" so it sits at the end of the\n# new last paragraph (right after \"Bhuratan\"), matching Word's normal\n# behavior of keeping _GoBack at the most recent edit location.\n$existing = $d.Bookmarks.Item(\"_GoBack\")\n$existing.Delete()\n\n# A bookmark collapsed exactly at the trailing edge of a paragraph's text\n# (the very last position before the paragraph mark) gets written out as\n# spanning the whole run, so a temporary trailing character is used to give\n# the collapsed point real content on both sides before it is removed again.\n$newParaRange = $newPara.Range\n$newParaRange.Collapse(0)\n$placeholderStart = $newParaRange.Start\n$newParaRange.InsertAfter(\"~\")\n\n$bookmarkRange = $d.Range($placeholderStart, $placeholderStart)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)\n$placeholderRange.Delete()\n"}
